$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: date serial, and B/C/D values (all zero)
$rows = @(
    @{ Row = 245; Date = 44319 },
    @{ Row = 246; Date = 44320 },
    @{ Row = 247; Date = 44321 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $prev = $r - 1

    # Copy formatting (style) from the cell directly above in column A
    $ws.Range("A$prev").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
